$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: NIM, NAMA, KODE_PROGRAM_STUDI, ANGKATAN (was NAMA, PROVINSI, KABUPATEN, KECAMATAN, KELURAHAN)
$ws.Range("A1").Value = "NIM"
$ws.Range("B1").Value = "NAMA"
$ws.Range("C1").Value = "KODE_PROGRAM_STUDI"
$ws.Range("D1").Value = "ANGKATAN"

# Remove the now-unused 5th column (KELURAHAN)
$ws.Range("E1").Clear()
